# Update "想去人数" (F column) counts on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1313
    $ws.Range("F3").Value = 1730
    $ws.Range("F5").Value = 6260
    $ws.Range("F6").Value = 98
    $ws.Range("F7").Value = 108
}
